# feat: add 2022-Q3 data
#
# The "2022-Q2" sheet (index 2) becomes the brand-new "2022-Q3" sheet,
# inserted right before the existing "2022-Q2" sheet (whose own data is
# untouched and simply slides one tab to the right, same as every other
# quarter sheet). The "总计" (totals) summary sheet gets its rows bumped
# down by one quarter and a new trailing row appended for "2021-Q2".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert a new "2022-Q3" worksheet, cloned from "2022-Q2" so it keeps
#    identical styling (bold/bordered header row, alignment, etc.), then
#    overwrite its data with the new quarter's numbers.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

$q3.Range("C2").Value = "工银全球精选股票（QDII）"

# Force these as literal text (matching the source file's inlineStr cells)
# instead of letting COM auto-coerce numeric-looking strings to numbers.
$q3.Range("D2:G2").NumberFormat = "@"
$q3.Range("D2").Value = "3.72"
$q3.Range("E2").Value = "93.69"
$q3.Range("F2").Value = "1.51"
$q3.Range("G2").Value = "0.0562"
$q3.Range("H2").Value = 10

# ---------------------------------------------------------------------
# 2) Update the "总计" summary sheet: every existing quarter row shifts
#    down (2022-Q2 -> row3, 2022-Q1 -> row4, ...) and a new row 7 is
#    appended for "2021-Q2".
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("A6:D6").Copy()
$total.Range("A7:D7").PasteSpecial(-4122)

$total.Range("B2").Value = "2022-Q3"
$total.Range("D2").Value = 0.06

$total.Range("B3").Value = "2022-Q2"
$total.Range("D3").Value = 0.08

$total.Range("B4").Value = "2022-Q1"
$total.Range("D4").Value = 0.08

$total.Range("B5").Value = "2021-Q4"
$total.Range("D5").Value = 0.09

$total.Range("B6").Value = "2021-Q3"
$total.Range("D6").Value = 0.08

$total.Range("A7").Value = 5
$total.Range("B7").Value = "2021-Q2"
$total.Range("C7").Value = 1
$total.Range("D7").Value = 0.08

# ---------------------------------------------------------------------
# 3) Keep the "2021-Q2" tab as the active/selected sheet, matching the
#    original workbook's tab selection.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q2").Activate()
